$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.742.99'
$ws.Range('E2').Value = '  -0.83%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.233.22'
$ws.Range('E3').Value = '  -2.71%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.645'
$ws.Range('E5').Value = '  +2.62%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '230.42'
$ws.Range('E6').Value = '  -0.37%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.20'
$ws.Range('E7').Value = '  +3.60%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.447'
$ws.Range('E9').Value = '  +4.91%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0959'
$ws.Range('E10').Value = '  +1.27%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.84'
$ws.Range('E11').Value = '  -1.81%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.50'
$ws.Range('E12').Value = '  +9.12%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').Value = '  +0.60%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.563.47'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.40'
$ws.Range('E15').Value = '  -2.44%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.11'
$ws.Range('E16').Value = '  +3.06%  '

$ws.Range('E17').Value = '  +0.82%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.235.37'
$ws.Range('E18').Value = '  -3.04%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.626.39'
$ws.Range('E19').Value = '  -0.87%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0981'
$ws.Range('E20').Value = '  +3.59%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.63'
$ws.Range('E21').Value = '  -1.54%  '

$ws.Range('E22').Value = '  -4.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '248.41'
$ws.Range('E23').Value = '  -2.48%  '

$ws.Range('E24').Value = '  -0.02%  '

$ws.Range('E25').Value = '  -5.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.39'
$ws.Range('E26').Value = '  +22.81%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.30'
$ws.Range('E27').Value = '  -2.95%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.83'
$ws.Range('E28').Value = '  -0.99%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.66'
$ws.Range('E29').Value = '  -0.35%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.76'
$ws.Range('E30').Value = '  +0.68%  '

$ws.Range('E31').Value = '  -2.13%  '

$ws.Range('E32').Value = '  -3.84%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.125'
$ws.Range('E33').Value = '  +2.58%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0695'
$ws.Range('E34').Value = '  +5.21%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.74'
$ws.Range('E35').Value = '  -1.66%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.88'
$ws.Range('E36').Value = '  -4.35%  '

$ws.Range('E37').Value = '  -0.48%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.37'
$ws.Range('E38').Value = '  -2.52%  '

$ws.Range('E39').Value = '  -6.10%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0259'
$ws.Range('E40').Value = '  +3.21%  '

$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000218'
$ws.Range('E42').Value = '  -2.96%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.23'
$ws.Range('E43').Value = '  -6.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.09'
$ws.Range('E44').Value = '  +0.21%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.94'
$ws.Range('E45').Value = '  -2.37%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0946'
$ws.Range('E46').Value = '  -2.50%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.35'
$ws.Range('E47').Value = '  -0.96%  '

$ws.Range('E48').Value = '  -2.90%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.32'
$ws.Range('E49').Value = '  +2.29%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.428.52'
$ws.Range('E50').Value = '  -3.71%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.75'
$ws.Range('E51').Value = '  +1.40%  '
